$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Body paragraph: append the GAP-closing carve-out sentence onto the
#    end of the existing disbursement-condition paragraph.
# ---------------------------------------------------------------------
$oldTail = "in accordance with the Approved Closing Statement."
$newTail = "in accordance with the Approved Closing Statement (or alternatively, the conditions set forth in paragraph 3(d) of paragraph C above for a GAP closing have been fully satisfied)."

$bodyFound = $d.Content.Find.Execute(
    $oldTail, $true, $false, $false, $false, $false,
    $true, 1, $false, $newTail, 2)

# ---------------------------------------------------------------------
# 2) Footer SAVEDATE field caches: "9-14-21" -> "3-9-22".
#    Walk every section's Footers collection (primary + first-page)
#    and replace the cached date text wherever that footer is actually
#    part of the document (HeaderFooter.Exists).
# ---------------------------------------------------------------------
$oldDate = "9-14-21"
$newDate = "3-9-22"

foreach ($sec in $d.Sections) {
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            $ftr.Range.Find.Execute(
                $oldDate, $true, $false, $false, $false, $false,
                $true, 1, $false, $newDate, 2) | Out-Null
        }
    }
}
